$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated Price values are plain decimals (e.g. "592.94") which Excel
# would otherwise auto-convert to a Number on assignment. Force those specific
# cells to Text format first so they stay strings, matching the other Price
# cells in the column that already look like text (e.g. "62.955.09").
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '62.955.09'
$ws.Range('D3').Value = '3.035.06'
$ws.Range('E3').Value = '  +1.13%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '592.94'
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('D6').Value = '153.50'
$ws.Range('E6').Value = '  +6.25%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.028.82'
$ws.Range('E8').Value = '  +0.93%  '
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('E10').Value = '  +7.81%  '
$ws.Range('E11').Value = '  +2.31%  '
$ws.Range('E12').Value = '  +0.25%  '
$ws.Range('E13').Value = '  +2.63%  '
$ws.Range('D14').Value = '35.56'
$ws.Range('E14').Value = '  +3.51%  '
$ws.Range('E15').Value = '  +1.93%  '
$ws.Range('D16').Value = '3.538.36'
$ws.Range('E16').Value = '  +1.23%  '
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('D18').Value = '62.941.87'
$ws.Range('E18').Value = '  +2.18%  '
$ws.Range('D19').Value = '3.034.47'
$ws.Range('E19').Value = '  +1.23%  '
$ws.Range('D20').Value = '452.35'
$ws.Range('E20').Value = '  -0.71%  '
$ws.Range('D21').Value = '14.30'
$ws.Range('E21').Value = '  +1.72%  '
$ws.Range('E22').Value = '  +0.63%  '
$ws.Range('E23').Value = '  +1.57%  '
$ws.Range('E24').Value = '  +0.99%  '
$ws.Range('E25').Value = '  +5.06%  '
$ws.Range('E26').Value = '  +5.99%  '
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('E29').Value = '  +7.87%  '
$ws.Range('E30').Value = '  +0.76%  '
$ws.Range('E31').Value = '  +7.80%  '
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('D33').Value = '27.55'
$ws.Range('E33').Value = '  +1.06%  '
$ws.Range('E34').Value = '  +2.85%  '
$ws.Range('D35').Value = '0.0₃0871'
$ws.Range('E35').Value = '  +5.94%  '
$ws.Range('E36').Value = '  +1.69%  '
$ws.Range('E37').Value = '  +2.73%  '
$ws.Range('D38').Value = '3.17'
$ws.Range('E38').Value = '  +9.56%  '
$ws.Range('E39').Value = '  +0.59%  '
$ws.Range('D40').Value = '50.56'
$ws.Range('E40').Value = '  +0.47%  '
$ws.Range('D41').Value = '9.11'
$ws.Range('E41').Value = '  -1.03%  '
$ws.Range('E42').Value = '  +3.49%  '
$ws.Range('E43').Value = '  +16.48%  '
$ws.Range('D44').Value = '42.42'
$ws.Range('E44').Value = '  +7.57%  '
$ws.Range('D45').Value = '397.11'
$ws.Range('E45').Value = '  -0.91%  '
$ws.Range('E46').Value = '  +1.74%  '
$ws.Range('D47').Value = '2.743.25'
$ws.Range('E47').Value = '  +0.78%  '
$ws.Range('D48').Value = '132.39'
$ws.Range('E48').Value = '  -0.60%  '
$ws.Range('E50').Value = '  +2.96%  '
$ws.Range('D51').Value = '24.27'
$ws.Range('E51').Value = '  +3.17%  '
